$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Results" sheet: insert a new column B (copy of column A), re-point the
# moved formulas, then overwrite a handful of the new column-B cells with
# the new "Miami" naming-convention strings.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Results")

# Duplicate column A into a freshly inserted column B (this carries over
# values, shared-string usage and cell styles in one shot).
$ws.Columns.Item(1).Copy()
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# The insert shifted the old B/C/D formula column to C/D/E; re-assign the
# difference formulas on the (now) E column so they stay linked as shared
# formulas (matching the original D-column grouping) and keep pointing at
# C/D instead of B/C.
$ws.Range("E2").Formula = "=C2-D2"
$ws.Range("E3:E13").Formula = "=C3-D3"
$ws.Range("E14:E45").Formula = "=C14-D14"
$ws.Range("E46:E69").Formula = "=C46-D46"
$ws.Range("E71").Formula = "=C71-D71"
$ws.Range("E72:E78").Formula = "=C72-D72"

# Overwrite specific column-B cells with the new Miami naming-convention
# strings (rest of column B keeps the value copied from column A).
$ws.Range("B16").Value = "PC_Local_to_Remote_MiamiM_full_bmp"
$ws.Range("B17").Value = "PC_Local_to_Remote_MiamiM_Double_bmp"
$ws.Range("B19").Value = "PC_Local_to_Remote_MiamiM_half_bmp"
$ws.Range("B20").Value = "PC_Local_to_Remote_MiamiM_half_png"
$ws.Range("B18").Value = "PC_Local_to_Remote_MiamiM_Double_png"
$ws.Range("B2").Value = "PC_Local_to_Remote_MiamiM_full_png"
$ws.Range("B15").Value = "PC_Local_to_Remote_MiamiM_full_png"
$ws.Range("B14").Value = "PC_Remote_to_Local_MiamiM_full_png"
$ws.Range("B3").Value = "PI_Remote_To_Local_MiamiM_full_PNG"

# ---------------------------------------------------------------------------
# View state: make "Results" the active sheet/tab, move each sheet's
# remembered selection.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Naming Conventions")
$ws1.Range("A2").Select()

$ws.Activate()
$ws.Range("A13").Select()
